# Locais.xlsx update — "Add files via upload"
#
# This adds a new "dia 1" (day 1) tag + "panda.png" icon to a handful of
# Edinburgh Harry-Potter-trail attractions (rows 27, 28, 34, 43), switches
# row 44's icon to harry-potter.png, renames the "dia 4" day-tag to "Bath"
# (shared by the Bath-related rows 47/53/54/55) and renumbers their "ordem"
# from 4 to 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the "dia 4" tag to "Bath" everywhere it is used -----------------
# (rows 47, 53, 54 and 55 all share this tag; writing the same new text to
# every cell that currently holds it lets the shared string be renamed in
# place instead of forking a duplicate entry.)
$ws.Range("H47").Value = "Bath"
$ws.Range("H53").Value = "Bath"
$ws.Range("H54").Value = "Bath"
$ws.Range("H55").Value = "Bath"

# Those same rows move from ordem 4 to ordem 2.
$ws.Range("J47").Value = 2
$ws.Range("J53").Value = 2
$ws.Range("J54").Value = 2
$ws.Range("J55").Value = 2

# --- Tag the Harry-Potter-trail stops as "dia 1" -----------------------------
$ws.Range("H27").Value = "dia 1"
$ws.Range("H28").Value = "dia 1"
$ws.Range("H34").Value = "dia 1"
$ws.Range("H43").Value = "dia 1"

# ordem = 1 for those new day-1 stops
$ws.Range("J27").Value = 1
$ws.Range("J28").Value = 1
$ws.Range("J34").Value = 1
$ws.Range("J43").Value = 1

# Swap their icon from the generic "touristic.png" to the new "panda.png"
$ws.Range("E27").Value = "panda.png"
$ws.Range("E28").Value = "panda.png"
$ws.Range("E34").Value = "panda.png"
$ws.Range("E43").Value = "panda.png"

# Row 44 (George Heriot's School) gets the harry-potter.png icon instead
$ws.Range("E44").Value = "harry-potter.png"

# --- Restore the workbook's on-screen selection ------------------------------
[void]$ws.Range("E43").Select()
